$wb = $excel.ActiveWorkbook

# --- Sheet "Vendas": rename D1 header and add new E1 column ---
$wsVendas = $wb.Worksheets.Item("Vendas")
$wsVendas.Range("D1").Value = "valor_venda"
$wsVendas.Range("E1").Value = "venda_id"

# Give the new header cell E1 the same (bold/centered) header style as the
# other header cells in row 1
$wsVendas.Range("D1").Copy()
$wsVendas.Range("E1").PasteSpecial(-4122)

# --- Sheet "Financeiro": fix column mismatch ---
$wsFin = $wb.Worksheets.Item("Financeiro")

# Rename header B1 from "categoria" to "tipo"
$wsFin.Range("B1").Value = "tipo"

# Move values from column E (old "tipo" data) into column B (which was empty)
$wsFin.Range("B2").Value = $wsFin.Range("E2").Value()
$wsFin.Range("B3").Value = $wsFin.Range("E3").Value()
$wsFin.Range("B4").Value = $wsFin.Range("E4").Value()

# Delete the now-redundant column E entirely
$wsFin.Range("E1:E4").Delete()
